$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values: Item changes from "Tepung Terigu"/"gr" to "Telur"/"pcs"
$ws.Range("A2").Value = "Telur"
$ws.Range("B2").Value = "2023-03-13"
$ws.Range("C2").Value = 16
$ws.Range("D2").Value = "pcs"
$ws.Range("E2").Value = 18000
$ws.Range("F2").Formula = "=E2/C2"

# Move the selection from E2 to A2
$ws.Range("A2").Select()

$wb.Save()
